# Regenerate save_data column G ("K", formerly Strike#) with new computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-11 (column G), replacing old Strike#-derived values.
$newValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 4
    6  = 0
    7  = 1
    8  = 2
    9  = 2
    10 = 2
    11 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
